{"js": "// The document has a paragraph containing a Word field:\n//   { m:'doc.html'.fromHTMLURI() }\n// stored as fldChar begin/instrText runs/fldChar end. The edit converts\n// that field into literal visible text runs \"{\", \"m\", \":\", \"'\", \"doc.html\",\n// \"'.fromHTMLURI()\", \"}\" (keeping the existing _GoBack bookmark in place),\n// i.e. the field delimiters become literal braces and the field code\n// becomes plain text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that holds the field (its fldChar/instrText content\n// is not exposed through the normal `.text` property, so we search the\n// paragraph's OOXML for the field markers instead of relying on `.text`).\nlet target = null;\nlet targetXml = null;\nfor (const p of paragraphs.items) {\n  p.__ooxml = p.getOoxml();\n}\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  const xml = p.__ooxml.value;\n  if (xml && xml.indexOf(\"fldChar\") !== -1 && xml.indexOf(\"instrText\") !== -1) {\n    target = p;\n    targetXml = xml;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the paragraph containing the field to rewrite.\");\n}\n\n// Preserve the paragraph's own rsid bookkeeping attributes exactly as they\n// already are in the document (drop any synthesized w14:paraId/textId that\n// getOoxml() adds on round-trip, those were not present originally).\nconst pOpenMatch = targetXml.match(/<w:p\\b([^>]*)>/);\nconst rawAttrs = pOpenMatch ? pOpenMatch[1] : \"\";\nconst attrRe = /(\\w+:\\w+)=\"([^\"]*)\"/g;\nlet pAttrs = \"\";\nlet am;\nwhile ((am = attrRe.exec(rawAttrs)) !== null) {\n  if (am[1].indexOf(\"w:rsid\") === 0) {\n    pAttrs += ` ${am[1]}=\"${am[2]}\"`;\n  }\n}\n\n// Rebuild the paragraph as literal text runs, keeping the existing\n// bookmark (_GoBack) in the same position between \"doc.html\" and\n// \"'.fromHTMLURI()\". The paragraph-level attributes (rsid info) are\n// preserved as-is from the original document.\nconst newParagraphXml =\n  '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"' +\n  pAttrs +\n  \">\" +\n  \"<w:r><w:t>{</w:t></w:r>\" +\n  \"<w:r><w:t>m</w:t></w:r>\" +\n  \"<w:r><w:t>:</w:t></w:r>\" +\n  \"<w:r><w:t>'</w:t></w:r>\" +\n  \"<w:r><w:t>doc.html</w:t></w:r>\" +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  \"<w:r><w:t>'.fromHTMLURI()</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  \"</w:p>\";\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  newParagraphXml +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nconst wholeRange = target.getRange(\"Whole\");\nwholeRange.insertOoxml(flatOpcXml, \"Replace\");\nawait context.sync();\n", "ps1": "# The document has a paragraph containing a Word field:\n#   { m:'doc.html'.fromHTMLURI() }\n# stored as fldChar begin/instrText runs/fldChar end. This script converts\n# that field into literal visible text runs \"{\", \"m\", \":\", \"'\", \"doc.html\",\n# \"'.fromHTMLURI()\", \"}\" (keeping the existing _GoBack bookmark in place),\n# i.e. the field delimiters become literal braces and the field code\n# becomes plain text.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that holds the field (search by Fields.Count rather\n# than a hardcoded index, so this keeps working if the document layout\n# shifts).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Fields.Count -gt 0) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the paragraph containing the field to rewrite.\"\n}\n\n$r = $target.Range\n\n# Pull the paragraph's own OOXML so we can keep its existing rsid\n# bookkeeping attributes intact (drop any synthesized w14:paraId/textId\n# that the round trip adds, those were not present originally).\n$fullXml = $r.XML()\n$bodyIdx = $fullXml.IndexOf(\"<w:body>\")\n$pIdx = $fullXml.IndexOf(\"<w:p\", $bodyIdx)\n$pEndIdx = $fullXml.IndexOf(\"</w:p>\", $pIdx) + 6\n$paraXml = $fullXml.Substring($pIdx, $pEndIdx - $pIdx)\n\n$openTagMatch = [regex]::Match($paraXml, '<w:p\\b([^>]*)>')\n$rawAttrs = $openTagMatch.Groups[1].Value\n\n$attrRe = [regex]'(\\w+:\\w+)=\"([^\"]*)\"'\n$pAttrs = \"\"\nforeach ($m in $attrRe.Matches($rawAttrs)) {\n    if ($m.Groups[1].Value.StartsWith(\"w:rsid\")) {\n        $pAttrs += \" \" + $m.Groups[1].Value + [char]61 + '\"' + $m.Groups[2].Value + '\"'\n    }\n}\n\n# Rebuild the paragraph as literal text runs, keeping the existing\n# bookmark (_GoBack) in the same position between \"doc.html\" and\n# \"'.fromHTMLURI()\".\n$newParagraphXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"' + $pAttrs + '>' + `\n    '<w:r><w:t>{</w:t></w:r>' + `\n    '<w:r><w:t>m</w:t></w:r>' + `\n    '<w:r><w:t>:</w:t></w:r>' + `\n    '<w:r><w:t>''</w:t></w:r>' + `\n    '<w:r><w:t>doc.html</w:t></w:r>' + `\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' + `\n    '<w:bookmarkEnd w:id=\"0\"/>' + `\n    '<w:r><w:t>''.fromHTMLURI()</w:t></w:r>' + `\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' + `\n    '</w:p>'\n\n$flatOpcXml = '<?xml version=\"1.0\" standalone=\"yes\"?>' + `\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n    '<pkg:xmlData>' + `\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n    '<w:body>' + $newParagraphXml + '</w:body>' + `\n    '</w:document>' + `\n    '</pkg:xmlData>' + `\n    '</pkg:part>' + `\n    '</pkg:package>'\n\n$r.InsertXML($flatOpcXml)\n"}
